$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) values were refreshed for both the
# "展览" sheet and the combined "全部类型" sheet, which mirrors the same
# rows. Update cell F2, F3, F4 on each.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1587
    $ws.Range("F3").Value = 107
    $ws.Range("F4").Value = 42
}
